$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 1070, shifting existing rows 1070:1146 down to 1071:1147.
$ws.Rows("1070:1070").Insert()

# Populate the newly inserted row 1070 with the new record values.
$ws.Range("A1070").Value = 9
$ws.Range("B1070").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1070").Value = "Metropolitana"
$ws.Range("D1070").Value = 45265
$ws.Range("E1070").Value = 13
$ws.Range("F1070").Value = 100112040
$ws.Range("G1070").Value = "Cilantro"
$ws.Range("H1070").Value = "Sin especificar"
$ws.Range("I1070").Value = "Primera"
$ws.Range("J1070").Value = 160
$ws.Range("K1070").Value = 20000
$ws.Range("L1070").Value = 22000
$ws.Range("M1070").Value = 21000
$ws.Range("N1070").Value = "$/docena de atados"
$ws.Range("O1070").Value = "Región Metropolitana"
$ws.Range("P1070").Value = 7000
$ws.Range("Q1070").Value = 3
$ws.Range("R1070").Value = "Hortaliza"

# Ensure the date cell keeps the expected date/time number format used by column D.
$ws.Range("D1070").NumberFormat = $ws.Range("D1071").NumberFormat
